$wb = $excel.ActiveWorkbook

# --- Service Contacts: move "delivery_organisation_path" column from R to D ---
# (before practitioner_key), shifting the intervening columns right by one.
$wsSC = $wb.Worksheets.Item("Service Contacts")
$wsSC.Columns("R").Cut()
$wsSC.Columns("D").Insert()

# --- View / selection changes on the affected sheets ---

$wsOrg = $wb.Worksheets.Item("Organisations")
$wsOrg.Range("H1:J3").Select()

$wsK10 = $wb.Worksheets.Item("K10+")
$wsK10.Range("F1:F5").Select()

$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Range("F1:F5").Select()

$wsSDQ = $wb.Worksheets.Item("SDQ")
$wsSDQ.Range("F4").Select()

# Leave "Service Contacts" as the active/selected sheet, with the new
# delivery_organisation_path column (D) selected.
$wsSC.Activate()
$wsSC.Columns("D").Select()
